# Updated cryptos list on Wed Oct 30 22:46:11 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's COM layer auto-converts plain numeric-looking strings (e.g. "26.20")
# assigned via Range.Value into actual Number cells, which both changes the
# stored cell type (t="n" instead of t="inlineStr"/shared-string "Text") and
# can silently drop meaningful trailing zeros ("26.20" -> 26.2). The source
# workbook stores every Price/Volume column as literal text, so force the
# "Text" number format before assigning, then restore the cell to the
# "Normal" style afterwards so no stray style index is left referenced on
# the cell (keeps styles.xml cell assignments identical to the original).
function Set-PriceText($row, $text) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-Volume($row, $text) {
    $ws.Cells.Item($row, 5).Value = $text
}

Set-PriceText 2  "72.408.56"
Set-Volume    2  "  -0.43%  "

Set-PriceText 3  "2.654.87"
Set-Volume    3  "  +0.61%  "

Set-Volume    4  "  +0.03%  "

Set-PriceText 5  "596.73"
Set-Volume    5  "  -1.64%  "

Set-PriceText 6  "175.21"
Set-Volume    6  "  -2.66%  "

Set-Volume    7  "  +0.03%  "

Set-Volume    8  "  -0.92%  "

Set-PriceText 9  "2.654.49"
Set-Volume    9  "  +0.63%  "

Set-Volume    10 "  -3.50%  "

Set-Volume    11 "  +1.96%  "

Set-Volume    12 "  +0.12%  "

Set-Volume    13 "  -1.15%  "

Set-PriceText 14 "3.138.07"
Set-Volume    14 "  +1.02%  "

Set-Volume    15 "  -2.71%  "

Set-PriceText 16 "72.373.61"
Set-Volume    16 "  -0.15%  "

Set-PriceText 17 "26.20"
Set-Volume    17 "  -2.43%  "

Set-PriceText 18 "2.654.33"
Set-Volume    18 "  +0.82%  "

Set-PriceText 19 "12.33"
Set-Volume    19 "  +5.15%  "

Set-PriceText 20 "370.66"
Set-Volume    20 "  -3.37%  "

Set-PriceText 21 "7.20"
Set-Volume    21 "  -9.75%  "

Set-Volume    22 "  -0.31%  "

Set-PriceText 23 "2.07"
Set-Volume    23 "  +0.68%  "

Set-PriceText 24 "71.99"
Set-Volume    24 "  -2.76%  "

Set-Volume    25 "  -0.13%  "

Set-PriceText 26 "4.32"
Set-Volume    26 "  -2.90%  "

Set-PriceText 27 "9.80"
Set-Volume    27 "  -2.36%  "

Set-PriceText 28 "2.797.69"
Set-Volume    28 "  +1.02%  "

Set-PriceText 29 "0.999"
Set-Volume    29 "  +0.01%  "

Set-Volume    30 "  +0.60%  "

Set-PriceText 31 "8.15"
Set-Volume    31 "  +0.32%  "

Set-PriceText 32 "495.15"
Set-Volume    32 "  -4.68%  "

Set-Volume    33 "  -2.99%  "

Set-Volume    34 "  -0.93%  "

Set-PriceText 35 "0.999"
Set-Volume    35 "  +0.05%  "

Set-PriceText 36 "162.09"
Set-Volume    36 "  -1.88%  "

Set-PriceText 37 "19.49"
Set-Volume    37 "  +0.19%  "

# Rows 38 and 39 swap positions (Kaspa <-> WhiteBITCoin), with the Kaspa
# price/volume also being refreshed to new values.
$ws.Cells.Item(38, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-PriceText 38 "18.94"
Set-Volume    38 "  -0.83%  "

$ws.Cells.Item(39, 2).Value = "Kaspa"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-PriceText 39 "0.111"
Set-Volume    39 "  -0.73%  "

Set-Volume    40 "  -2.99%  "

Set-Volume    41 "  -5.41%  "

Set-Volume    42 "  -0.04%  "

Set-Volume    43 "  -3.44%  "

Set-Volume    44 "  -0.20%  "

Set-Volume    45 "  -1.06%  "

Set-PriceText 46 "155.70"
Set-Volume    46 "  +3.23%  "

Set-PriceText 47 "39.22"
Set-Volume    47 "  -0.63%  "

Set-Volume    48 "  +0.53%  "

Set-Volume    49 "  +1.66%  "

Set-Volume    50 "  +1.11%  "

Set-Volume    51 "  -1.40%  "
